# Actualización automática 2025-09-12 16:35:09
#
# "VENTAS POR GRUPO" : L19/M19 gain sales for client LLERENA CONDO SANDRA
#   MARISOL (PIEDRA SINTERIZADA / PORCELANATO groups), and the row-35
#   "N de 33" completion counters for those two columns bump by one.
# "VENTA MENSUAL"    : F19 (septiembre) for the same client, and the F35
#   column total.
# "CUMPLIMIENTO MENSUAL" : the PIEDRA SINTERIZADA / PORCELANATO / TOTAL
#   rows (11, 12, 15) get their VENTA (D), POR CUMPLIR (E) and
#   CUMPLIMIENTO (F) figures recomputed from the new sales.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 19 = ALMEIDA CUATIN JHONATHANN CARLOS / LLERENA CONDO SANDRA MARISOL
$wsGrupo.Cells.Item(19, 12).Value = 3224.04   # L19 (PIEDRA SINTERIZADA)
$wsGrupo.Cells.Item(19, 13).Value = 853.42    # M19 (PORCELANATO)

# Row 35 = completion counters "N de 33" for each group column
$wsGrupo.Cells.Item(35, 12).Value = "1 de 33" # L35
$wsGrupo.Cells.Item(35, 13).Value = "6 de 33" # M35

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL"
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Row 19 = ALMEIDA CUATIN JHONATHANN CARLOS / LLERENA CONDO SANDRA MARISOL
$wsMensual.Cells.Item(19, 6).Value = 4077.46  # F19 (septiembre)

# Row 35 = column totals
$wsMensual.Cells.Item(35, 6).Value = 13699.06 # F35 (septiembre total)

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 11 = PIEDRA SINTERIZADA
$wsCumpl.Cells.Item(11, 4).Value = 3224.04            # D11 VENTA
$wsCumpl.Cells.Item(11, 5).Value = -301.8154181472601 # E11 POR CUMPLIR
$wsCumpl.Cells.Item(11, 6).Value = 1.103282759313421  # F11 CUMPLIMIENTO

# Row 12 = PORCELANATO
$wsCumpl.Cells.Item(12, 4).Value = 8642.42            # D12 VENTA
$wsCumpl.Cells.Item(12, 5).Value = 13791.3353751766   # E12 POR CUMPLIR
$wsCumpl.Cells.Item(12, 6).Value = 0.3852417865607562 # F12 CUMPLIMIENTO

# Row 15 = TOTAL
$wsCumpl.Cells.Item(15, 4).Value = 13956.29           # D15 VENTA
$wsCumpl.Cells.Item(15, 5).Value = 24786.72881339593  # E15 POR CUMPLIR
$wsCumpl.Cells.Item(15, 6).Value = 0.3602272209922481 # F15 CUMPLIMIENTO

# Column E ("POR CUMPLIR") widened by one unit, matching Excel's
# column-autofit side effect of the updated values. (COM's ColumnWidth
# snaps to the workbook's pixel grid on save, so 23.165 -- not 24 -- is
# the input that round-trips to a stored width of exactly 24.)
$wsCumpl.Columns.Item(5).ColumnWidth = 23.165
